$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))
$lo2 = $ws.ListObjects.Item(1)
Write-Host "after append-resize:"
for ($i = 1; $i -le $lo2.ListColumns.Count; $i++) {
    Write-Host $i $lo2.ListColumns.Item($i).Name
}

# Now set values for O and P (the newly appended columns) to be authentic_source_id/name? No -
# Actually per target, these should end up mid-table. Let's instead try inserting columns now
# physically at J:K, see if the table column objects (ids 15/16) "follow" the cells they were bound to (O/P -> still O/P after insertion point at J, since insert is before O/P, so O/P's *content* shifts right to Q/R? no wait J:K insert shifts everything from J onwards right by 2)
$ws.Range("J1:K1").EntireColumn.Insert()
Write-Host "after column insert:"
$lo3 = $ws.ListObjects.Item(1)
Write-Host $lo3.Range.Address
for ($i = 1; $i -le $lo3.ListColumns.Count; $i++) {
    Write-Host $i $lo3.ListColumns.Item($i).Name
}
